# =====================================================================
# Commit: "feat: add 2022-Q3 data"
#
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计" (pushing
#    "2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q1" one slot later),
#    seeded with a copy of the "2022-Q2" sheet layout/formatting (header
#    row + index-column style) so the new sheet matches the workbook's
#    established per-quarter fund-holding table format.
# 2) Fill it with the 20 fund rows for 2022-Q3.
# 3) Update the "总计" (summary) sheet: add the 2022-Q3 row at the top
#    of the data (20 holdings / 4.82 亿元) and shift every other quarter
#    row down by one.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---- Step 1: insert + seed the new "2022-Q3" worksheet ----
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"
$ws = $wb.Worksheets.Item("2022-Q3")

# Seed layout/styles from "2022-Q2" (header row style + index-column style)
# then trim away the extra rows it does not need (2022-Q2 has 32 funds,
# 2022-Q3 only has 20).
$q2Sheet.Range("A1:H33").Copy($ws.Range("A1"))
$ws.Range("A1").Clear()
$ws.Range("A22:H33").Clear()

# Columns B (fund code) and D:G (scale/position/ratio/value) hold
# numeric-looking text in the source data - force text format so COM
# does not silently coerce e.g. "513060" or "44.00" into numbers.
$ws.Range("B2:B21").NumberFormat = "@"
$ws.Range("D2:G19").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "513060"
$ws.Range("C2").Value = "博时恒生医疗保健ETF（QDII）"
$ws.Range("D2").Value = "44.00"
$ws.Range("E2").Value = "99.48"
$ws.Range("F2").Value = "3.19"
$ws.Range("G2").Value = "1.4036"
$ws.Range("H2").Value = 8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "008954"
$ws.Range("C3").Value = "安信价值回报三年持有期混合A"
$ws.Range("D3").Value = "29.83"
$ws.Range("E3").Value = "90.84"
$ws.Range("F3").Value = "3.82"
$ws.Range("G3").Value = "1.1395"
$ws.Range("H3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "012892"
$ws.Range("C4").Value = "安信优质企业三年持有混合A"
$ws.Range("D4").Value = "14.77"
$ws.Range("E4").Value = "90.16"
$ws.Range("F4").Value = "4.63"
$ws.Range("G4").Value = "0.6839"
$ws.Range("H4").Value = 10

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "002387"
$ws.Range("C5").Value = "工银沪港深股票A"
$ws.Range("D5").Value = "13.37"
$ws.Range("E5").Value = "83.30"
$ws.Range("F5").Value = "3.96"
$ws.Range("G5").Value = "0.5295"
$ws.Range("H5").Value = 8

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "012893"
$ws.Range("C6").Value = "安信优质企业三年持有混合C"
$ws.Range("D6").Value = "6.68"
$ws.Range("E6").Value = "90.16"
$ws.Range("F6").Value = "4.63"
$ws.Range("G6").Value = "0.3093"
$ws.Range("H6").Value = 10

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "005197"
$ws.Range("C7").Value = "工银瑞信沪港深精选灵活配置混合A"
$ws.Range("D7").Value = "4.01"
$ws.Range("E7").Value = "94.39"
$ws.Range("F7").Value = "4.22"
$ws.Range("G7").Value = "0.1692"
$ws.Range("H7").Value = 9

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "005504"
$ws.Range("C8").Value = "汇添富沪港深大盘价值混合A"
$ws.Range("D8").Value = "3.14"
$ws.Range("E8").Value = "91.24"
$ws.Range("F8").Value = "4.50"
$ws.Range("G8").Value = "0.1413"
$ws.Range("H8").Value = 8

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "513700"
$ws.Range("C9").Value = "鹏华中证港股通医药卫生综合ETF"
$ws.Range("D9").Value = "2.96"
$ws.Range("E9").Value = "94.14"
$ws.Range("F9").Value = "2.66"
$ws.Range("G9").Value = "0.0787"
$ws.Range("H9").Value = 9

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "470888"
$ws.Range("C10").Value = "汇添富香港优势精选混合（QDII）"
$ws.Range("D10").Value = "1.63"
$ws.Range("E10").Value = "78.50"
$ws.Range("F10").Value = "4.44"
$ws.Range("G10").Value = "0.0724"
$ws.Range("H10").Value = 6

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "159892"
$ws.Range("C11").Value = "华夏恒生香港上市生物科技ETF（QDII）"
$ws.Range("D11").Value = "1.58"
$ws.Range("E11").Value = "99.13"
$ws.Range("F11").Value = "3.51"
$ws.Range("G11").Value = "0.0555"
$ws.Range("H11").Value = 8

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "513120"
$ws.Range("C12").Value = "广发中证香港创新药（QDII-ETF）"
$ws.Range("D12").Value = "1.09"
$ws.Range("E12").Value = "98.58"
$ws.Range("F12").Value = "4.96"
$ws.Range("G12").Value = "0.0541"
$ws.Range("H12").Value = 7

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "513280"
$ws.Range("C13").Value = "汇添富恒生香港上市生物科技ETF（QDII）"
$ws.Range("D13").Value = "1.51"
$ws.Range("E13").Value = "100.14"
$ws.Range("F13").Value = "3.55"
$ws.Range("G13").Value = "0.0536"
$ws.Range("H13").Value = 8

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "005198"
$ws.Range("C14").Value = "工银瑞信沪港深精选灵活配置混合C"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "94.39"
$ws.Range("F14").Value = "4.22"
$ws.Range("G14").Value = "0.0426"
$ws.Range("H14").Value = 9

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "007512"
$ws.Range("C15").Value = "工银沪港深股票C"
$ws.Range("D15").Value = "0.86"
$ws.Range("E15").Value = "83.30"
$ws.Range("F15").Value = "3.96"
$ws.Range("G15").Value = "0.0341"
$ws.Range("H15").Value = 8

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "513200"
$ws.Range("C16").Value = "易方达中证港股通医药卫生综合ETF"
$ws.Range("D16").Value = "0.77"
$ws.Range("E16").Value = "95.67"
$ws.Range("F16").Value = "2.78"
$ws.Range("G16").Value = "0.0214"
$ws.Range("H16").Value = 9

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "159776"
$ws.Range("C17").Value = "银华中证港股通医药卫生综合ETF"
$ws.Range("D17").Value = "0.52"
$ws.Range("E17").Value = "92.74"
$ws.Range("F17").Value = "2.62"
$ws.Range("G17").Value = "0.0136"
$ws.Range("H17").Value = 9

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "159718"
$ws.Range("C18").Value = "平安中证港股通医药卫生综合ETF"
$ws.Range("D18").Value = "0.53"
$ws.Range("E18").Value = "90.14"
$ws.Range("F18").Value = "2.55"
$ws.Range("G18").Value = "0.0135"
$ws.Range("H18").Value = 9

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "010667"
$ws.Range("C19").Value = "安信价值回报三年持有期混合C"
$ws.Range("D19").Value = "0.03"
$ws.Range("E19").Value = "90.84"
$ws.Range("F19").Value = "3.82"
$ws.Range("G19").Value = "0.0011"
$ws.Range("H19").Value = 10

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "015118"
$ws.Range("C20").Value = "汇添富沪港深大盘价值混合C"
$ws.Range("D20").Value = "0.00"
$ws.Range("E20").Value = "91.24"
$ws.Range("F20").Value = "4.50"
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 8

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "015119"
$ws.Range("C21").Value = "汇添富沪港深大盘价值混合D"
$ws.Range("D21").Value = "0.00"
$ws.Range("E21").Value = "91.24"
$ws.Range("F21").Value = "4.50"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 8

# ---- Step 2: update the "总计" summary sheet ----
$totalWs = $wb.Worksheets.Item("总计")

# Row 7 is brand new (the table previously only spanned rows 1-6) - copy
# the existing index-column cell style (s="2") down from A2 before
# writing into it, so A7 is formatted the same as A2:A6.
$totalWs.Range("A2").Copy($totalWs.Range("A7"))

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 20
$totalWs.Range("D2").Value = 4.82

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q2"
$totalWs.Range("C3").Value = 32
$totalWs.Range("D3").Value = 8.06

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2022-Q1"
$totalWs.Range("C4").Value = 34
$totalWs.Range("D4").Value = 10.44

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q4"
$totalWs.Range("C5").Value = 29
$totalWs.Range("D5").Value = 11.02

$totalWs.Range("A6").Value = 4
$totalWs.Range("B6").Value = "2021-Q3"
$totalWs.Range("C6").Value = 62
$totalWs.Range("D6").Value = 24.62

$totalWs.Range("A7").Value = 5
$totalWs.Range("B7").Value = "2021-Q1"
$totalWs.Range("C7").Value = 22
$totalWs.Range("D7").Value = 13.13
